$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new value looks like a plain number so that Excel
# keeps storing them as text (matching the original inline-string data),
# instead of silently converting them to numeric values.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated cell values
$ws.Range('D2').Value = '27.188.72'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '1.631.50'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '216.61'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').Value = '0.518'
$ws.Range('E6').Value = '  +1.50%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('D9').Value = '0.0625'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('D10').Value = '20.35'
$ws.Range('E10').Value = '  +1.96%  '
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '1.630.29'
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').Value = '27.189.41'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = '64.83'
$ws.Range('E16').Value = '  -3.86%  '
$ws.Range('D17').Value = '0.0₃0734'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').Value = '215.66'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('E23').Value = '  -0.98%  '
$ws.Range('D24').Value = '148.10'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('E26').Value = '  -1.95%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('D31').Value = '3.39'
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('E32').Value = '  -0.54%  '
$ws.Range('D33').Value = '1.311.95'
$ws.Range('E33').Value = '  +3.46%  '
$ws.Range('D34').Value = '1.57'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').Value = '2.49'
$ws.Range('E35').Value = '  +1.49%  '
$ws.Range('E36').Value = '  -1.43%  '
$ws.Range('D37').Value = '0.851'
$ws.Range('E37').Value = '  +1.31%  '
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('E40').Value = '  +1.66%  '
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('D42').Value = '63.80'
$ws.Range('E42').Value = '  +2.26%  '
$ws.Range('D43').Value = '1.770.31'
$ws.Range('D44').Value = '5.21'
$ws.Range('E44').Value = '  -3.59%  '
$ws.Range('D45').Value = '90.76'
$ws.Range('E45').Value = '  -1.07%  '
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').Value = '0.804'
$ws.Range('E48').Value = '  +20.53%  '
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').Value = '7.57'
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.0951'
$ws.Range('E51').Value = '  -2.41%  '
